$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 318.35715
$ws.Range("I33").Value = 331.30768
$ws.Range("K33").Value = 331.30768
$ws.Range("M33").Value = -102.30768
$ws.Range("H40").Value = 41673476
$ws.Range("I40").Value = 7500
$ws.Range("K40").Value = 7500
$ws.Range("M40").Value = -7325
$ws.Range("H54").Value = 20146
$ws.Range("I54").Value = 14166.667
$ws.Range("K54").Value = 14166.667
$ws.Range("M54").Value = -13680.667
$ws.Range("H62").Value = 7058.9
$ws.Range("I62").Value = 6882.3335
$ws.Range("J62").Value = 7323.75
$ws.Range("K62").Value = 6882.3335
$ws.Range("L62").Value = 7323.75
$ws.Range("M62").Value = -6258.3335
$ws.Range("N62").Value = -8571.75
$ws.Range("H65").Value = 7058.9
$ws.Range("I65").Value = 6882.3335
$ws.Range("J65").Value = 7323.75
$ws.Range("K65").Value = 34411.6675
$ws.Range("L65").Value = 36618.75
$ws.Range("M65").Value = -31291.6675
$ws.Range("N65").Value = -42858.75
$ws.Range("H70").Value = 2034705.1
$ws.Range("J70").Value = 2952.5
$ws.Range("L70").Value = 8857.5
$ws.Range("N70").Value = -9397.5
$ws.Range("H73").Value = 2034705.1
$ws.Range("J73").Value = 2952.5
$ws.Range("L73").Value = 8857.5
$ws.Range("N73").Value = -10729.5
$ws.Range("H113").Value = 7565.231
$ws.Range("I113").Value = 6699.8335
$ws.Range("J113").Value = 8307
$ws.Range("K113").Value = 6699.8335
$ws.Range("L113").Value = 8307
$ws.Range("M113").Value = -3445.8335
$ws.Range("N113").Value = -14815
$ws.Range("H129").Value = 4551.636
$ws.Range("I129").Value = 890.3333
$ws.Range("K129").Value = 2670.9999
$ws.Range("M129").Value = 2329.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6575.8867
$ws.Range("I32").Value = 6074.961
$ws.Range("K32").Value = 6074.961
$ws.Range("M32").Value = -5787.961
$ws.Range("H74").Value = 2132.9333
$ws.Range("I74").Value = 1273.091
$ws.Range("J74").Value = 4497.5
$ws.Range("K74").Value = 1273.091
$ws.Range("L74").Value = 4497.5
$ws.Range("M74").Value = -399.0909999999999
$ws.Range("N74").Value = -6245.5
$ws.Range("H77").Value = 2132.9333
$ws.Range("I77").Value = 1273.091
$ws.Range("J77").Value = 4497.5
$ws.Range("K77").Value = 6365.455
$ws.Range("L77").Value = 22487.5
$ws.Range("M77").Value = -1997.455
$ws.Range("N77").Value = -31223.5
$ws.Range("H101").Value = 94999.5
$ws.Range("J101").Value = 94999.5
$ws.Range("L101").Value = 94999.5
$ws.Range("N101").Value = -101489.5
$ws.Range("H110").Value = 3797.9666
$ws.Range("I110").Value = 3431.1428
$ws.Range("K110").Value = 3431.1428
$ws.Range("M110").Value = -1386.1428
$ws.Range("H132").Value = 2633583.5
$ws.Range("I132").Value = 1899.2424
$ws.Range("K132").Value = 5697.7272
$ws.Range("M132").Value = -3167.7272

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4847.1113
$ws.Range("I107").Value = 5730.2856
$ws.Range("K107").Value = 5730.2856
$ws.Range("M107").Value = -3810.2856

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 81999.5
$ws.Range("J18").Value = 81999.5
$ws.Range("L18").Value = 81999.5
$ws.Range("N18").Value = -82459.5
$ws.Range("H58").Value = 3902.6956
$ws.Range("I58").Value = 3450.9333
$ws.Range("J58").Value = 4749.75
$ws.Range("K58").Value = 3450.9333
$ws.Range("L58").Value = 4749.75
$ws.Range("M58").Value = -3247.9333
$ws.Range("N58").Value = -5155.75
$ws.Range("H99").Value = 18685.25
$ws.Range("I99").Value = 9455.5
$ws.Range("K99").Value = 9455.5
$ws.Range("M99").Value = -7957.5
$ws.Range("H126").Value = 18685.25
$ws.Range("I126").Value = 9455.5
$ws.Range("K126").Value = 28366.5
$ws.Range("M126").Value = -25896.5
$ws.Range("H132").Value = 2144.3572
$ws.Range("I132").Value = 1812.2
$ws.Range("K132").Value = 5436.6
$ws.Range("M132").Value = -2906.6
$ws.Range("H136").Value = 3902.6956
$ws.Range("I136").Value = 3450.9333
$ws.Range("J136").Value = 4749.75
$ws.Range("K136").Value = 10352.7999
$ws.Range("L136").Value = 14249.25
$ws.Range("M136").Value = -7802.7999
$ws.Range("N136").Value = -19349.25
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 984.7143
$ws.Range("J23").Value = 984.7143
$ws.Range("L23").Value = 2954.1429
$ws.Range("N23").Value = -3424.1429
$ws.Range("H109").Value = 4167.1177
$ws.Range("I109").Value = 1250.5714
$ws.Range("K109").Value = 3751.7142
$ws.Range("M109").Value = -2711.7142
$ws.Range("H132").Value = 1716.5758
$ws.Range("I132").Value = 914.2
$ws.Range("J132").Value = 2065.4348
$ws.Range("K132").Value = 8227.800000000001
$ws.Range("L132").Value = 18588.9132
$ws.Range("M132").Value = -5697.800000000001
$ws.Range("N132").Value = -23648.9132

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2215.1667
$ws.Range("I102").Value = 2280.182
$ws.Range("K102").Value = 2280.182
$ws.Range("M102").Value = -658.1819999999998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7513.85
$ws.Range("I7").Value = 7237.722
$ws.Range("K7").Value = 7237.722
$ws.Range("M7").Value = -7125.722
$ws.Range("H22").Value = 7509.6523
$ws.Range("I22").Value = 14862
$ws.Range("J22").Value = 3588.4
$ws.Range("K22").Value = 14862
$ws.Range("L22").Value = 3588.4
$ws.Range("M22").Value = -14567
$ws.Range("N22").Value = -4178.4
$ws.Range("H27").Value = 7509.6523
$ws.Range("I27").Value = 14862
$ws.Range("J27").Value = 3588.4
$ws.Range("K27").Value = 14862
$ws.Range("L27").Value = 3588.4
$ws.Range("M27").Value = -14755
$ws.Range("N27").Value = -3802.4
$ws.Range("H55").Value = 1125.3125
$ws.Range("J55").Value = 1643
$ws.Range("L55").Value = 1643
$ws.Range("N55").Value = -1989
$ws.Range("H100").Value = 16687025
$ws.Range("I100").Value = 2616.25
$ws.Range("J100").Value = 35754920
$ws.Range("K100").Value = 2616.25
$ws.Range("L100").Value = 35754920
$ws.Range("M100").Value = -2075.25
$ws.Range("N100").Value = -35756002
$ws.Range("H122").Value = 4437.551
$ws.Range("I122").Value = 3401.1462
$ws.Range("K122").Value = 10203.4386
$ws.Range("M122").Value = -7753.438600000001
$ws.Range("H126").Value = 7513.85
$ws.Range("I126").Value = 7237.722
$ws.Range("K126").Value = 21713.166
$ws.Range("M126").Value = -19243.166
$ws.Range("H136").Value = 2244.6667
$ws.Range("I136").Value = 2093.7
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 6281.099999999999
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -3731.099999999999
$ws.Range("N136").Value = -14098.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 14166.333
$ws.Range("I28").Value = 6995
$ws.Range("K28").Value = 6995
$ws.Range("M28").Value = -6647
$ws.Range("H113").Value = 476.44
$ws.Range("I113").Value = 351.8125
$ws.Range("K113").Value = 1055.4375
$ws.Range("M113").Value = 1114.5625
$ws.Range("H126").Value = 4029.25
$ws.Range("I126").Value = 4691.3076
$ws.Range("K126").Value = 14073.9228
$ws.Range("M126").Value = -11603.9228
$ws.Range("H128").Value = 99999
$ws.Range("J128").Value = 99999
$ws.Range("L128").Value = 99999
$ws.Range("N128").Value = -109959
$ws.Range("H137").Value = 98888
$ws.Range("J137").Value = 98888
$ws.Range("L137").Value = 98888
$ws.Range("N137").Value = -109088
$ws.Range("H141").Value = 86527.39999999999
$ws.Range("J141").Value = 86527.39999999999
$ws.Range("L141").Value = 86527.39999999999
$ws.Range("N141").Value = -96887.39999999999

Write-Host "Applied all Ragnarok_Profits market-data updates"